# Update "想去人数" (interest count) figures in the three data sheets
# that track the Hefei comic-convention listings ("展览", "演出", "全部类型").
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 513
$ws1.Range("F4").Value  = 187
$ws1.Range("F6").Value  = 492
$ws1.Range("F8").Value  = 105
$ws1.Range("F9").Value  = 38
$ws1.Range("F10").Value = 6508
$ws1.Range("F11").Value = 220
$ws1.Range("F12").Value = 352
$ws1.Range("F13").Value = 2672
$ws1.Range("F14").Value = 162
$ws1.Range("F15").Value = 275
$ws1.Range("F16").Value = 253
$ws1.Range("F17").Value = 507

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 11

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 11
$ws4.Range("F5").Value  = 513
$ws4.Range("F6").Value  = 187
$ws4.Range("F8").Value  = 492
$ws4.Range("F10").Value = 105
$ws4.Range("F11").Value = 38
$ws4.Range("F13").Value = 6508
$ws4.Range("F15").Value = 220
$ws4.Range("F16").Value = 352
$ws4.Range("F17").Value = 2672
$ws4.Range("F18").Value = 162
$ws4.Range("F19").Value = 275
$ws4.Range("F20").Value = 253
$ws4.Range("F21").Value = 507
